# Apply scheduled market-data refresh to Sheets (Leve profit calculations)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2377
$ws.Range("I19").Value = 1445
$ws.Range("K19").Value = 1445
$ws.Range("M19").Value = -1270
$ws.Range("H33").Value = 382.5
$ws.Range("I33").Value = 241
$ws.Range("K33").Value = 241
$ws.Range("M33").Value = -12
$ws.Range("H62").Value = 9532267
$ws.Range("J62").Value = 6966.3335
$ws.Range("L62").Value = 6966.3335
$ws.Range("N62").Value = -8214.333500000001
$ws.Range("H65").Value = 9532267
$ws.Range("J65").Value = 6966.3335
$ws.Range("L65").Value = 34831.6675
$ws.Range("N65").Value = -41071.6675
$ws.Range("H86").Value = 71431160
$ws.Range("I86").Value = 100002540
$ws.Range("J86").Value = 2699.5
$ws.Range("K86").Value = 100002540
$ws.Range("L86").Value = 2699.5
$ws.Range("M86").Value = -100001417
$ws.Range("N86").Value = -4945.5
$ws.Range("H89").Value = 71431160
$ws.Range("I89").Value = 100002540
$ws.Range("J89").Value = 2699.5
$ws.Range("K89").Value = 500012700
$ws.Range("L89").Value = 13497.5
$ws.Range("M89").Value = -500007084
$ws.Range("N89").Value = -24729.5
$ws.Range("H98").Value = 2064.5925
$ws.Range("I98").Value = 2132.077
$ws.Range("K98").Value = 2132.077
$ws.Range("M98").Value = -634.0770000000002
$ws.Range("H107").Value = 1118.5
$ws.Range("I107").Value = 1101.9131
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1101.9131
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 818.0869
$ws.Range("N107").Value = -5340
$ws.Range("H111").Value = 3549.6
$ws.Range("I111").Value = 2916
$ws.Range("K111").Value = 8748
$ws.Range("M111").Value = -5681
$ws.Range("H122").Value = 2064.5925
$ws.Range("I122").Value = 2132.077
$ws.Range("K122").Value = 6396.231000000001
$ws.Range("M122").Value = -3946.231000000001
$ws.Range("H137").Value = 12754.952
$ws.Range("I137").Value = 7546.4614
$ws.Range("J137").Value = 21218.75
$ws.Range("K137").Value = 22639.3842
$ws.Range("L137").Value = 63656.25
$ws.Range("M137").Value = -20089.3842
$ws.Range("N137").Value = -68756.25
$ws.Range("H138").Value = 3881.8914
$ws.Range("J138").Value = 6111.36
$ws.Range("L138").Value = 18334.08
$ws.Range("N138").Value = -28614.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3145.6667
$ws.Range("I45").Value = 2773.7778
$ws.Range("J45").Value = 4261.3335
$ws.Range("K45").Value = 2773.7778
$ws.Range("L45").Value = 4261.3335
$ws.Range("M45").Value = -2396.7778
$ws.Range("N45").Value = -5015.3335
$ws.Range("H74").Value = 3500.8667
$ws.Range("I74").Value = 3736.5454
$ws.Range("K74").Value = 3736.5454
$ws.Range("M74").Value = -2862.5454
$ws.Range("H77").Value = 3500.8667
$ws.Range("I77").Value = 3736.5454
$ws.Range("K77").Value = 18682.727
$ws.Range("M77").Value = -14314.727
$ws.Range("H88").Value = 7226.727
$ws.Range("I88").Value = 2066.6667
$ws.Range("J88").Value = 9161.75
$ws.Range("K88").Value = 2066.6667
$ws.Range("L88").Value = 9161.75
$ws.Range("M88").Value = -1660.6667
$ws.Range("N88").Value = -9973.75
$ws.Range("H91").Value = 7226.727
$ws.Range("I91").Value = 2066.6667
$ws.Range("J91").Value = 9161.75
$ws.Range("K91").Value = 2066.6667
$ws.Range("L91").Value = 9161.75
$ws.Range("M91").Value = -662.6667000000002
$ws.Range("N91").Value = -11969.75
$ws.Range("H97").Value = 1831.6
$ws.Range("I97").Value = 1119.3334
$ws.Range("K97").Value = 1119.3334
$ws.Range("M97").Value = -623.3334
$ws.Range("H110").Value = 3714.6216
$ws.Range("I110").Value = 1180.1177
$ws.Range("J110").Value = 5868.95
$ws.Range("K110").Value = 1180.1177
$ws.Range("L110").Value = 5868.95
$ws.Range("M110").Value = 864.8823
$ws.Range("N110").Value = -9958.950000000001
$ws.Range("H132").Value = 706273.4399999999
$ws.Range("I132").Value = 825291.7
$ws.Range("K132").Value = 2475875.1
$ws.Range("M132").Value = -2473345.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 15164006
$ws.Range("J80").Value = 22239528
$ws.Range("L80").Value = 22239528
$ws.Range("N80").Value = -22241524
$ws.Range("H83").Value = 15164006
$ws.Range("J83").Value = 22239528
$ws.Range("L83").Value = 111197640
$ws.Range("N83").Value = -111207624
$ws.Range("H88").Value = 17250
$ws.Range("J88").Value = 17250
$ws.Range("L88").Value = 17250
$ws.Range("N88").Value = -18062
$ws.Range("H91").Value = 17250
$ws.Range("J91").Value = 17250
$ws.Range("L91").Value = 17250
$ws.Range("N91").Value = -20058
$ws.Range("H99").Value = 7018.427
$ws.Range("I99").Value = 6889.2666
$ws.Range("J99").Value = 7212.1665
$ws.Range("K99").Value = 6889.2666
$ws.Range("L99").Value = 7212.1665
$ws.Range("M99").Value = -5391.2666
$ws.Range("N99").Value = -10208.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6286.933
$ws.Range("I86").Value = 3888.05
$ws.Range("K86").Value = 3888.05
$ws.Range("M86").Value = -2765.05
$ws.Range("H89").Value = 6286.933
$ws.Range("I89").Value = 3888.05
$ws.Range("K89").Value = 19440.25
$ws.Range("M89").Value = -13824.25
$ws.Range("H99").Value = 7411507.5
$ws.Range("I99").Value = 11114887
$ws.Range("K99").Value = 11114887
$ws.Range("M99").Value = -11113389
$ws.Range("H105").Value = 59830080
$ws.Range("I105").Value = 89743920
$ws.Range("K105").Value = 89743920
$ws.Range("M105").Value = -89742173
$ws.Range("H107").Value = 849.5454999999999
$ws.Range("I107").Value = 359.83334
$ws.Range("K107").Value = 359.83334
$ws.Range("M107").Value = 1560.16666
$ws.Range("H126").Value = 7411507.5
$ws.Range("I126").Value = 11114887
$ws.Range("K126").Value = 33344661
$ws.Range("M126").Value = -33342191
$ws.Range("H132").Value = 7395.2383
$ws.Range("I132").Value = 6529.4707
$ws.Range("J132").Value = 11074.75
$ws.Range("K132").Value = 19588.4121
$ws.Range("L132").Value = 33224.25
$ws.Range("M132").Value = -17058.4121
$ws.Range("N132").Value = -38284.25
$ws.Range("H134").Value = 58833120
$ws.Range("I134").Value = 100009430
$ws.Range("J134").Value = 9820.571
$ws.Range("K134").Value = 300028290
$ws.Range("L134").Value = 29461.713
$ws.Range("M134").Value = -300025755
$ws.Range("N134").Value = -34531.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 386.66666
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H52").Value = 5499.5
$ws.Range("J52").Value = 5499.5
$ws.Range("L52").Value = 16498.5
$ws.Range("N52").Value = -17030.5
$ws.Range("H131").Value = 37684140
$ws.Range("I131").Value = 41028496
$ws.Range("K131").Value = 123085488
$ws.Range("M131").Value = -123080448
$ws.Range("H132").Value = 2389.8096
$ws.Range("I132").Value = 1467.2858
$ws.Range("K132").Value = 13205.5722
$ws.Range("M132").Value = -10675.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 116161.445
$ws.Range("I122").Value = 128212.875
$ws.Range("J122").Value = 19750
$ws.Range("K122").Value = 384638.625
$ws.Range("L122").Value = 59250
$ws.Range("M122").Value = -382188.625
$ws.Range("N122").Value = -64150

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3321.5386
$ws.Range("I40").Value = 3016.7144
$ws.Range("K40").Value = 3016.7144
$ws.Range("M40").Value = -2880.7144
$ws.Range("H61").Value = 6225.387
$ws.Range("I61").Value = 5641.778
$ws.Range("K61").Value = 5641.778
$ws.Range("M61").Value = -5439.778
$ws.Range("H113").Value = 6225.387
$ws.Range("I113").Value = 5641.778
$ws.Range("K113").Value = 5641.778
$ws.Range("M113").Value = -3471.778
